$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Fix the "No." column style on rows 19, 22 and 25 (s="2" -> s="5", the
#    style already used for the shaded "No." cells elsewhere, e.g. A6).
# ---------------------------------------------------------------------------
$ws.Range("A6").Copy()
$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("A22").PasteSpecial(-4122)
$ws.Range("A25").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. Append four new fixture rows (26-29).
#    Formats are seeded from existing rows that already carry the exact
#    combination of styles needed, then the values are written on top so no
#    new styles are minted (except the one genuinely new style used below).
# ---------------------------------------------------------------------------

# --- Row 26: No./Date styled like row 2 & 7, Team/Venue/Time styled like row 3
$ws.Range("A2:B2").Copy()
$ws.Range("A26:B26").PasteSpecial(-4122)
$ws.Range("B7").Copy()
$ws.Range("B26").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("C3:G3").Copy()
$ws.Range("C26:G26").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("D26").PasteSpecial(-4122)

$ws.Range("A26").Value2 = 25
$ws.Range("B26").Value2 = 43631
$ws.Range("C26").Value = "Westridge Warriors"
$ws.Range("D26").Value = "Maharashtra Sports Club"
$ws.Range("E26").Value = "Russell Creek G3 P or R"
$ws.Range("F26").Value = "07.30 AM"
$ws.Range("G26").Value = "11.30 AM"

# --- Row 27: everything styled like row 2, except the date cell (row 7 style)
$ws.Range("A2:G2").Copy()
$ws.Range("A27:G27").PasteSpecial(-4122)
$ws.Range("B7").Copy()
$ws.Range("B27").PasteSpecial(-4122)

$ws.Range("A27").Value2 = 26
$ws.Range("B27").Value2 = 43631
$ws.Range("C27").Value = "LazzyLADS"
$ws.Range("D27").Value = "Mustangs Blue"
$ws.Range("E27").Value = "Russell Creek G5 S or T"
$ws.Range("F27").Value = "11.30 AM"
$ws.Range("G27").Value = "03.30 PM"

# --- Row 28: No. styled like row 3 (shaded), rest like row 2, date like row 7
$ws.Range("A2:G2").Copy()
$ws.Range("A28:G28").PasteSpecial(-4122)
$ws.Range("A3").Copy()
$ws.Range("A28").PasteSpecial(-4122)
$ws.Range("B7").Copy()
$ws.Range("B28").PasteSpecial(-4122)

$ws.Range("A28").Value2 = 27
$ws.Range("B28").Value2 = 43631
$ws.Range("C28").Value = "423 Spartans"
$ws.Range("D28").Value = "Plano Titans"
$ws.Range("E28").Value = "Frisco Independence"
$ws.Range("F28").Value = "07.30 AM"
$ws.Range("G28").Value = "11.30 AM"

# --- Row 29: like row 2 (A,D,F,G), date like row 7, Team1 cell (C) like row 3.
#     E29 is intentionally left at the sheet's default style so the
#     HorizontalAlignment assignment below mints a genuinely new style.
$ws.Range("A2:D2").Copy()
$ws.Range("A29:D29").PasteSpecial(-4122)
$ws.Range("F2:G2").Copy()
$ws.Range("F29:G29").PasteSpecial(-4122)
$ws.Range("B7").Copy()
$ws.Range("B29").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C29").PasteSpecial(-4122)

$ws.Range("A29").Value2 = 28
$ws.Range("B29").Value2 = 43632
$ws.Range("C29").Value = "Cruisers"
$ws.Range("D29").Value = "Shadow Warriors"
$ws.Range("G29").Value = "07.30 PM"

# F29 stays blank (matches source fixture sheet); E29 gets a brand new,
# left-aligned general style and the new venue string.
$ws.Range("E29").HorizontalAlignment = -4131
$ws.Range("E29").Value = "Russell Creek G4 L or M"

# ---------------------------------------------------------------------------
# 3. Update the active selection to the newly added F29 cell.
# ---------------------------------------------------------------------------
$ws.Range("F29").Select()
